# demo2-excel to two dim string array
# Adds 10 new rows (rows 4-13) of invalid-credential test data to the
# "invalidCredentialTest" sheet (sheet1.xml / ActiveSheet), using a
# two-dimensional string array, username "Abi" and a sequence of
# passwords admin123, admin124 ... admin132, all flagged with the
# existing "Invalid credentials" message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$username = "Abi"
$errorMsg = "Invalid credentials"

# Two dimensional string array: [password, errorMessage] per extra row.
$data = @(
    @("admin123", $errorMsg),
    @("admin124", $errorMsg),
    @("admin125", $errorMsg),
    @("admin126", $errorMsg),
    @("admin127", $errorMsg),
    @("admin128", $errorMsg),
    @("admin129", $errorMsg),
    @("admin130", $errorMsg),
    @("admin131", $errorMsg),
    @("admin132", $errorMsg)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $username
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}

$lastRow = $startRow + $data.Length - 1
$ws.Range("A4:C$lastRow").Select()
